# Update countries & provincias Spain
# Refresh the COVID figures for a handful of countries (new daily numbers),
# update the "last updated" timestamp, then re-sort the country table by
# "Casos totales" (column B) descending - same as the live data refresh that
# produced the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last updated" banner -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 23:42"

# --- 2. Update per-country figures (row numbers as currently laid out) -
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6705507
$ws.Range("C4").Value = 28906
$ws.Range("D4").Value = 3966181
$ws.Range("E4").Value = 2540859
$ws.Range("G4").Value = 339
$ws.Range("H4").Value = 198467

# Brasil (row 6)
$ws.Range("B6").Value = 4330455
$ws.Range("C6").Value = 14597
$ws.Range("E6").Value = 645409
$ws.Range("G6").Value = 351
$ws.Range("H6").Value = 131625

# Ecuador (row 32)
$ws.Range("B32").Value = 118594
$ws.Range("C32").Value = 2143
$ws.Range("D32").Value = 97063
$ws.Range("E32").Value = 10628
$ws.Range("G32").Value = 39
$ws.Range("H32").Value = 10903

# Guinea (row 96)
$ws.Range("B96").Value = 10045
$ws.Range("C96").Value = 25
$ws.Range("D96").Value = 9292
$ws.Range("E96").Value = 690

# Haiti (row 104)
$ws.Range("B104").Value = 8493
$ws.Range("C104").Value = 15
$ws.Range("E104").Value = 2154

# Cabo Verde (row 121)
$ws.Range("B121").Value = 4813
$ws.Range("C121").Value = 102
$ws.Range("D121").Value = 4119
$ws.Range("E121").Value = 650
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 44

# Cuba (row 122)
$ws.Range("B122").Value = 4684
$ws.Range("C122").Value = 31
$ws.Range("D122").Value = 3930
$ws.Range("E122").Value = 646

# Ruanda (row 124)
$ws.Range("B124").Value = 4591
$ws.Range("C124").Value = 26
$ws.Range("D124").Value = 2556
$ws.Range("E124").Value = 2013
$ws.Range("H124").Value = 22

# Angola (row 132)
$ws.Range("B132").Value = 3388
$ws.Range("C132").Value = 53
$ws.Range("D132").Value = 1301
$ws.Range("E132").Value = 1953
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 134

# Togo (row 157)
$ws.Range("B157").Value = 1572
$ws.Range("C157").Value = 17
$ws.Range("D157").Value = 1190
$ws.Range("E157").Value = 345

# --- 3. Re-sort the whole country table by Casos totales (desc) --------
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2)
